$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = -7.378800000000004
$ws.Range("B9").Value = 5.349399999999996
$ws.Range("D12").Value = -6.954800000000001
$ws.Range("B13").Value = 6.323599999999999
$ws.Range("D14").Value = -7.881600000000001
$ws.Range("B16").Value = 5.900000000000002
$ws.Range("B18").Value = 7.118099999999994
$ws.Range("D19").Value = -7.789800000000001
$ws.Range("B20").Value = 8.699200000000001
$ws.Range("B26").Value = 5.435400000000007
$ws.Range("D26").Value = -8.878399999999997
$ws.Range("B27").Value = 5.627500000000003
$ws.Range("D27").Value = -8.8466
$ws.Range("B29").Value = 5.035799999999999
$ws.Range("D29").Value = -7.288399999999995
$ws.Range("B35").Value = 8.286800000000005
$ws.Range("B36").Value = 8.942100000000005
$ws.Range("D37").Value = -7.660599999999999
$ws.Range("D38").Value = -8.008199999999999
$ws.Range("B45").Value = 4.938500000000006
$ws.Range("D47").Value = -7.317000000000002
$ws.Range("D51").Value = -8.150899999999998
$ws.Range("D52").Value = -7.4544
$ws.Range("B55").Value = 7.115199999999995
$ws.Range("D55").Value = -7.744600000000003
$ws.Range("B57").Value = 5.148499999999997
$ws.Range("B69").Value = 5.153299999999998
$ws.Range("D69").Value = -7.142199999999995
$ws.Range("D70").Value = -7.568100000000002
$ws.Range("B76").Value = 5.078599999999999
$ws.Range("D76").Value = -8.103699999999995
$ws.Range("B78").Value = 10.0577
$ws.Range("D81").Value = -7.773400000000001
$ws.Range("B82").Value = 6.2626
$ws.Range("B83").Value = 5.354800000000001
$ws.Range("D83").Value = -9.149299999999993
$ws.Range("B93").Value = 5.765
$ws.Range("D94").Value = -7.163399999999998
$ws.Range("B97").Value = 6.232699999999998
$ws.Range("D100").Value = -8.423500000000001
$ws.Range("D102").Value = -7.783299999999997
